# ---------------------------------------------------------------------------
# Sync attendance_reports: add a new MICROBIOLOGY session (24/11/2025) on the
# 'Attendance' log and roll the per-student stats on 'Summary' forward for the
# 24 students who attended it.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")
$attendance = $wb.Worksheets.Item("Attendance")

# --- 1) Append the 24 new attendance-log rows (688 -> 712) ------------------
$attendance.Range("A689").Value = "'221031"
$attendance.Range("B689").Value = "امنيه عبدالله عبد اللطيف محمد"
$attendance.Range("C689").Value = "Year 2"
$attendance.Range("D689").Value = "C1"
$attendance.Range("E689").Value = "221031@med.asu.edu.eg"
$attendance.Range("F689").Value = "MICROBIOLOGY"
$attendance.Range("G689").Value = "'1"
$attendance.Range("H689").Value = "MICROBIOLOGY"
$attendance.Range("I689").Value = "24/11/2025"
$attendance.Range("J689").Value = "09:21:14"
$attendance.Range("K689").Value = "C1"

$attendance.Range("A690").Value = "'221584"
$attendance.Range("B690").Value = "عزه بنت محمد بن عوض الصمداني"
$attendance.Range("C690").Value = "Year 2"
$attendance.Range("D690").Value = "C1"
$attendance.Range("E690").Value = "221584@med.asu.edu.eg"
$attendance.Range("F690").Value = "MICROBIOLOGY"
$attendance.Range("G690").Value = "'1"
$attendance.Range("H690").Value = "MICROBIOLOGY"
$attendance.Range("I690").Value = "24/11/2025"
$attendance.Range("J690").Value = "09:21:33"
$attendance.Range("K690").Value = "C1"

$attendance.Range("A691").Value = "'221307"
$attendance.Range("B691").Value = "دعاء عاصم على العوض"
$attendance.Range("C691").Value = "Year 2"
$attendance.Range("D691").Value = "C1"
$attendance.Range("E691").Value = "221307@med.asu.edu.eg"
$attendance.Range("F691").Value = "MICROBIOLOGY"
$attendance.Range("G691").Value = "'1"
$attendance.Range("H691").Value = "MICROBIOLOGY"
$attendance.Range("I691").Value = "24/11/2025"
$attendance.Range("J691").Value = "09:21:41"
$attendance.Range("K691").Value = "C1"

$attendance.Range("A692").Value = "'221818"
$attendance.Range("B692").Value = "جيهان محارب الشيخ الكيلاني"
$attendance.Range("C692").Value = "Year 2"
$attendance.Range("D692").Value = "C1"
$attendance.Range("E692").Value = "221818@med.asu.edu.eg"
$attendance.Range("F692").Value = "MICROBIOLOGY"
$attendance.Range("G692").Value = "'1"
$attendance.Range("H692").Value = "MICROBIOLOGY"
$attendance.Range("I692").Value = "24/11/2025"
$attendance.Range("J692").Value = "09:21:52"
$attendance.Range("K692").Value = "C1"

$attendance.Range("A693").Value = "'221810"
$attendance.Range("B693").Value = "رهان محارب الشيخ الكيلاني"
$attendance.Range("C693").Value = "Year 2"
$attendance.Range("D693").Value = "C1"
$attendance.Range("E693").Value = "221810@med.asu.edu.eg"
$attendance.Range("F693").Value = "MICROBIOLOGY"
$attendance.Range("G693").Value = "'1"
$attendance.Range("H693").Value = "MICROBIOLOGY"
$attendance.Range("I693").Value = "24/11/2025"
$attendance.Range("J693").Value = "09:22:09"
$attendance.Range("K693").Value = "C1"

$attendance.Range("A694").Value = "'221838"
$attendance.Range("B694").Value = "الاء سيد احمد احمد ابن ادريس"
$attendance.Range("C694").Value = "Year 2"
$attendance.Range("D694").Value = "C1"
$attendance.Range("E694").Value = "221838@med.asu.edu.eg"
$attendance.Range("F694").Value = "MICROBIOLOGY"
$attendance.Range("G694").Value = "'1"
$attendance.Range("H694").Value = "MICROBIOLOGY"
$attendance.Range("I694").Value = "24/11/2025"
$attendance.Range("J694").Value = "09:22:20"
$attendance.Range("K694").Value = "C1"

$attendance.Range("A695").Value = "'222035"
$attendance.Range("B695").Value = "علا عبد الوهاب خليل محمود"
$attendance.Range("C695").Value = "Year 2"
$attendance.Range("D695").Value = "C1"
$attendance.Range("E695").Value = "222035@med.asu.edu.eg"
$attendance.Range("F695").Value = "MICROBIOLOGY"
$attendance.Range("G695").Value = "'1"
$attendance.Range("H695").Value = "MICROBIOLOGY"
$attendance.Range("I695").Value = "24/11/2025"
$attendance.Range("J695").Value = "09:22:30"
$attendance.Range("K695").Value = "C1"

$attendance.Range("A696").Value = "'221319"
$attendance.Range("B696").Value = "روان صلاح طاهر الوهباني"
$attendance.Range("C696").Value = "Year 2"
$attendance.Range("D696").Value = "C1"
$attendance.Range("E696").Value = "221319@med.asu.edu.eg"
$attendance.Range("F696").Value = "MICROBIOLOGY"
$attendance.Range("G696").Value = "'1"
$attendance.Range("H696").Value = "MICROBIOLOGY"
$attendance.Range("I696").Value = "24/11/2025"
$attendance.Range("J696").Value = "09:22:39"
$attendance.Range("K696").Value = "C1"

$attendance.Range("A697").Value = "'210967"
$attendance.Range("B697").Value = "ملاك كمال اسماعيل ابو جلاله"
$attendance.Range("C697").Value = "Year 2"
$attendance.Range("D697").Value = "C1"
$attendance.Range("E697").Value = "210967@med.asu.edu.eg"
$attendance.Range("F697").Value = "MICROBIOLOGY"
$attendance.Range("G697").Value = "'1"
$attendance.Range("H697").Value = "MICROBIOLOGY"
$attendance.Range("I697").Value = "24/11/2025"
$attendance.Range("J697").Value = "09:23:01"
$attendance.Range("K697").Value = "C1"

$attendance.Range("A698").Value = "'222113"
$attendance.Range("B698").Value = "ابرار محمد عبد الله عبد الحميد"
$attendance.Range("C698").Value = "Year 2"
$attendance.Range("D698").Value = "C1"
$attendance.Range("E698").Value = "222113@med.asu.edu.eg"
$attendance.Range("F698").Value = "MICROBIOLOGY"
$attendance.Range("G698").Value = "'1"
$attendance.Range("H698").Value = "MICROBIOLOGY"
$attendance.Range("I698").Value = "24/11/2025"
$attendance.Range("J698").Value = "09:23:10"
$attendance.Range("K698").Value = "C1"

$attendance.Range("A699").Value = "'221675"
$attendance.Range("B699").Value = "ساره بنت سعيد بن عثمان الكناني"
$attendance.Range("C699").Value = "Year 2"
$attendance.Range("D699").Value = "C1"
$attendance.Range("E699").Value = "221675@med.asu.edu.eg"
$attendance.Range("F699").Value = "MICROBIOLOGY"
$attendance.Range("G699").Value = "'1"
$attendance.Range("H699").Value = "MICROBIOLOGY"
$attendance.Range("I699").Value = "24/11/2025"
$attendance.Range("J699").Value = "09:23:20"
$attendance.Range("K699").Value = "C1"

$attendance.Range("A700").Value = "'221527"
$attendance.Range("B700").Value = "شهد محمد عبدالرحمن ادريس"
$attendance.Range("C700").Value = "Year 2"
$attendance.Range("D700").Value = "C1"
$attendance.Range("E700").Value = "221527@med.asu.edu.eg"
$attendance.Range("F700").Value = "MICROBIOLOGY"
$attendance.Range("G700").Value = "'1"
$attendance.Range("H700").Value = "MICROBIOLOGY"
$attendance.Range("I700").Value = "24/11/2025"
$attendance.Range("J700").Value = "09:23:44"
$attendance.Range("K700").Value = "C1"

$attendance.Range("A701").Value = "'211704"
$attendance.Range("B701").Value = "الياس احمد بكردان"
$attendance.Range("C701").Value = "Year 2"
$attendance.Range("D701").Value = "C1"
$attendance.Range("E701").Value = "211704@med.asu.edu.eg"
$attendance.Range("F701").Value = "MICROBIOLOGY"
$attendance.Range("G701").Value = "'1"
$attendance.Range("H701").Value = "MICROBIOLOGY"
$attendance.Range("I701").Value = "24/11/2025"
$attendance.Range("J701").Value = "09:23:51"
$attendance.Range("K701").Value = "C1"

$attendance.Range("A702").Value = "'221324"
$attendance.Range("B702").Value = "مجد ذوقان خليل قيشاوي"
$attendance.Range("C702").Value = "Year 2"
$attendance.Range("D702").Value = "C1"
$attendance.Range("E702").Value = "221324@med.asu.edu.eg"
$attendance.Range("F702").Value = "MICROBIOLOGY"
$attendance.Range("G702").Value = "'1"
$attendance.Range("H702").Value = "MICROBIOLOGY"
$attendance.Range("I702").Value = "24/11/2025"
$attendance.Range("J702").Value = "09:23:57"
$attendance.Range("K702").Value = "C1"

$attendance.Range("A703").Value = "'221566"
$attendance.Range("B703").Value = "مصطفى سامى محمد عبد الله"
$attendance.Range("C703").Value = "Year 2"
$attendance.Range("D703").Value = "C1"
$attendance.Range("E703").Value = "221566@med.asu.edu.eg"
$attendance.Range("F703").Value = "MICROBIOLOGY"
$attendance.Range("G703").Value = "'1"
$attendance.Range("H703").Value = "MICROBIOLOGY"
$attendance.Range("I703").Value = "24/11/2025"
$attendance.Range("J703").Value = "09:24:09"
$attendance.Range("K703").Value = "C1"

$attendance.Range("A704").Value = "'221569"
$attendance.Range("B704").Value = "هبه جعفر محمد شوكت"
$attendance.Range("C704").Value = "Year 2"
$attendance.Range("D704").Value = "C1"
$attendance.Range("E704").Value = "221569@med.asu.edu.eg"
$attendance.Range("F704").Value = "MICROBIOLOGY"
$attendance.Range("G704").Value = "'1"
$attendance.Range("H704").Value = "MICROBIOLOGY"
$attendance.Range("I704").Value = "24/11/2025"
$attendance.Range("J704").Value = "09:24:13"
$attendance.Range("K704").Value = "C1"

$attendance.Range("A705").Value = "'212163"
$attendance.Range("B705").Value = "رقيه احمد عبد الله"
$attendance.Range("C705").Value = "Year 2"
$attendance.Range("D705").Value = "C1"
$attendance.Range("E705").Value = "212163@med.asu.edu.eg"
$attendance.Range("F705").Value = "MICROBIOLOGY"
$attendance.Range("G705").Value = "'1"
$attendance.Range("H705").Value = "MICROBIOLOGY"
$attendance.Range("I705").Value = "24/11/2025"
$attendance.Range("J705").Value = "09:24:52"
$attendance.Range("K705").Value = "C1"

$attendance.Range("A706").Value = "'221996"
$attendance.Range("B706").Value = "نياقوط فال توت دوير"
$attendance.Range("C706").Value = "Year 2"
$attendance.Range("D706").Value = "C1"
$attendance.Range("E706").Value = "221996@med.asu.edu.eg"
$attendance.Range("F706").Value = "MICROBIOLOGY"
$attendance.Range("G706").Value = "'1"
$attendance.Range("H706").Value = "MICROBIOLOGY"
$attendance.Range("I706").Value = "24/11/2025"
$attendance.Range("J706").Value = "09:25:08"
$attendance.Range("K706").Value = "C1"

$attendance.Range("A707").Value = "'221558"
$attendance.Range("B707").Value = "محمد عادل عوض باحاج"
$attendance.Range("C707").Value = "Year 2"
$attendance.Range("D707").Value = "C1"
$attendance.Range("E707").Value = "221558@med.asu.edu.eg"
$attendance.Range("F707").Value = "MICROBIOLOGY"
$attendance.Range("G707").Value = "'1"
$attendance.Range("H707").Value = "MICROBIOLOGY"
$attendance.Range("I707").Value = "24/11/2025"
$attendance.Range("J707").Value = "09:25:20"
$attendance.Range("K707").Value = "C1"

$attendance.Range("A708").Value = "'221494"
$attendance.Range("B708").Value = "حسن الصادق مصطفى الحاج"
$attendance.Range("C708").Value = "Year 2"
$attendance.Range("D708").Value = "C1"
$attendance.Range("E708").Value = "221494@med.asu.edu.eg"
$attendance.Range("F708").Value = "MICROBIOLOGY"
$attendance.Range("G708").Value = "'1"
$attendance.Range("H708").Value = "MICROBIOLOGY"
$attendance.Range("I708").Value = "24/11/2025"
$attendance.Range("J708").Value = "09:25:24"
$attendance.Range("K708").Value = "C1"

$attendance.Range("A709").Value = "'221459"
$attendance.Range("B709").Value = "محمد الطيب محمد زين"
$attendance.Range("C709").Value = "Year 2"
$attendance.Range("D709").Value = "C1"
$attendance.Range("E709").Value = "221459@med.asu.edu.eg"
$attendance.Range("F709").Value = "MICROBIOLOGY"
$attendance.Range("G709").Value = "'1"
$attendance.Range("H709").Value = "MICROBIOLOGY"
$attendance.Range("I709").Value = "24/11/2025"
$attendance.Range("J709").Value = "09:25:36"
$attendance.Range("K709").Value = "C1"

$attendance.Range("A710").Value = "'221536"
$attendance.Range("B710").Value = "عبده دفع الله سليمان كوكو"
$attendance.Range("C710").Value = "Year 2"
$attendance.Range("D710").Value = "C1"
$attendance.Range("E710").Value = "221536@med.asu.edu.eg"
$attendance.Range("F710").Value = "MICROBIOLOGY"
$attendance.Range("G710").Value = "'1"
$attendance.Range("H710").Value = "MICROBIOLOGY"
$attendance.Range("I710").Value = "24/11/2025"
$attendance.Range("J710").Value = "09:25:44"
$attendance.Range("K710").Value = "C1"

$attendance.Range("A711").Value = "'221522"
$attendance.Range("B711").Value = "آدم محمد احمد البديرات"
$attendance.Range("C711").Value = "Year 2"
$attendance.Range("D711").Value = "C1"
$attendance.Range("E711").Value = "221522@med.asu.edu.eg"
$attendance.Range("F711").Value = "MICROBIOLOGY"
$attendance.Range("G711").Value = "'1"
$attendance.Range("H711").Value = "MICROBIOLOGY"
$attendance.Range("I711").Value = "24/11/2025"
$attendance.Range("J711").Value = "09:25:51"
$attendance.Range("K711").Value = "C1"

$attendance.Range("A712").Value = "'221546"
$attendance.Range("B712").Value = "محمدزين ابوبكر محمد زين احمد"
$attendance.Range("C712").Value = "Year 2"
$attendance.Range("D712").Value = "C1"
$attendance.Range("E712").Value = "221546@med.asu.edu.eg"
$attendance.Range("F712").Value = "MICROBIOLOGY"
$attendance.Range("G712").Value = "'1"
$attendance.Range("H712").Value = "MICROBIOLOGY"
$attendance.Range("I712").Value = "24/11/2025"
$attendance.Range("J712").Value = "09:26:09"
$attendance.Range("K712").Value = "C1"

# --- 2) Re-stamp the autofilter + defined name for the grown range ----------
$attendance.AutoFilterMode = $false
$attendance.Range("A1:K712").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Attendance!_FilterDatabase") {
        $n.RefersTo = "='Attendance'!`$A`$1:`$K`$712"
    }
}

# --- 3) Roll forward the per-student MICROBIOLOGY stats on 'Summary' --------
# Row 23
$summary.Range("G23").Value = "'13.8%"
$summary.Range("I23").Value = 19
$summary.Range("N23").Value = 4
$summary.Range("O23").Value = 8
$summary.Range("AG23").Value = 1

# Row 38
$summary.Range("G38").Value = "'24.1%"
$summary.Range("I38").Value = 16
$summary.Range("N38").Value = 7
$summary.Range("O38").Value = 5
$summary.Range("AG38").Value = 1

# Row 55
$summary.Range("G55").Value = "'13.8%"
$summary.Range("I55").Value = 19
$summary.Range("N55").Value = 4
$summary.Range("O55").Value = 8
$summary.Range("AG55").Value = 1

# Row 100
$summary.Range("G100").Value = "'17.2%"
$summary.Range("I100").Value = 18
$summary.Range("N100").Value = 5
$summary.Range("O100").Value = 7
$summary.Range("AG100").Value = 1

# Row 108
$summary.Range("G108").Value = "'17.2%"
$summary.Range("I108").Value = 18
$summary.Range("N108").Value = 5
$summary.Range("O108").Value = 7
$summary.Range("AG108").Value = 1

# Row 110
$summary.Range("G110").Value = "'17.2%"
$summary.Range("I110").Value = 18
$summary.Range("N110").Value = 5
$summary.Range("O110").Value = 7
$summary.Range("AG110").Value = 1

# Row 111
$summary.Range("G111").Value = "'17.2%"
$summary.Range("I111").Value = 18
$summary.Range("N111").Value = 5
$summary.Range("O111").Value = 7
$summary.Range("AG111").Value = 1

# Row 130
$summary.Range("G130").Value = "'24.1%"
$summary.Range("I130").Value = 16
$summary.Range("N130").Value = 7
$summary.Range("O130").Value = 5
$summary.Range("AG130").Value = 1

# Row 134
$summary.Range("F134").Value = "Moderate Risk"
$summary.Range("F134").Interior.Color = 8239615
$summary.Range("G134").Value = "'27.6%"
$summary.Range("I134").Value = 15
$summary.Range("N134").Value = 8
$summary.Range("O134").Value = 4
$summary.Range("AG134").Value = 1

# Row 143
$summary.Range("G143").Value = "'17.2%"
$summary.Range("I143").Value = 18
$summary.Range("N143").Value = 5
$summary.Range("O143").Value = 7
$summary.Range("AG143").Value = 1

# Row 145
$summary.Range("F145").Value = "High Risk"
$summary.Range("F145").Interior.Color = 8158463
$summary.Range("G145").Value = "'20.7%"
$summary.Range("I145").Value = 17
$summary.Range("N145").Value = 6
$summary.Range("O145").Value = 6
$summary.Range("AG145").Value = 1

# Row 146
$summary.Range("F146").Value = "High Risk"
$summary.Range("F146").Interior.Color = 8158463
$summary.Range("G146").Value = "'20.7%"
$summary.Range("I146").Value = 17
$summary.Range("N146").Value = 6
$summary.Range("O146").Value = 6
$summary.Range("AG146").Value = 1

# Row 148
$summary.Range("G148").Value = "'17.2%"
$summary.Range("I148").Value = 18
$summary.Range("N148").Value = 5
$summary.Range("O148").Value = 7
$summary.Range("AG148").Value = 1

# Row 150
$summary.Range("G150").Value = "'10.3%"
$summary.Range("I150").Value = 20
$summary.Range("N150").Value = 3
$summary.Range("O150").Value = 9
$summary.Range("AG150").Value = 1

# Row 151
$summary.Range("G151").Value = "'10.3%"
$summary.Range("I151").Value = 20
$summary.Range("N151").Value = 3
$summary.Range("O151").Value = 9
$summary.Range("AG151").Value = 1

# Row 152
$summary.Range("F152").Value = "High Risk"
$summary.Range("F152").Interior.Color = 8158463
$summary.Range("G152").Value = "'20.7%"
$summary.Range("I152").Value = 17
$summary.Range("N152").Value = 6
$summary.Range("O152").Value = 6
$summary.Range("AG152").Value = 1

# Row 157
$summary.Range("G157").Value = "'10.3%"
$summary.Range("I157").Value = 20
$summary.Range("N157").Value = 3
$summary.Range("O157").Value = 9
$summary.Range("AG157").Value = 1

# Row 175
$summary.Range("F175").Value = "High Risk"
$summary.Range("F175").Interior.Color = 8158463
$summary.Range("G175").Value = "'20.7%"
$summary.Range("I175").Value = 17
$summary.Range("N175").Value = 6
$summary.Range("O175").Value = 6
$summary.Range("AG175").Value = 1

# Row 193
$summary.Range("G193").Value = "'10.3%"
$summary.Range("I193").Value = 20
$summary.Range("N193").Value = 3
$summary.Range("O193").Value = 9
$summary.Range("AG193").Value = 1

# Row 194
$summary.Range("G194").Value = "'10.3%"
$summary.Range("I194").Value = 20
$summary.Range("N194").Value = 3
$summary.Range("O194").Value = 9
$summary.Range("AG194").Value = 1

# Row 198
$summary.Range("G198").Value = "'13.8%"
$summary.Range("I198").Value = 19
$summary.Range("N198").Value = 4
$summary.Range("O198").Value = 8
$summary.Range("AG198").Value = 1

# Row 228
$summary.Range("G228").Value = "'17.2%"
$summary.Range("I228").Value = 18
$summary.Range("N228").Value = 5
$summary.Range("O228").Value = 7
$summary.Range("AG228").Value = 1

# Row 241
$summary.Range("F241").Value = "Moderate Risk"
$summary.Range("F241").Interior.Color = 8239615
$summary.Range("G241").Value = "'27.6%"
$summary.Range("I241").Value = 15
$summary.Range("N241").Value = 8
$summary.Range("O241").Value = 4
$summary.Range("AG241").Value = 1

# Row 252
$summary.Range("F252").Value = "High Risk"
$summary.Range("F252").Interior.Color = 8158463
$summary.Range("G252").Value = "'20.7%"
$summary.Range("I252").Value = 17
$summary.Range("N252").Value = 6
$summary.Range("O252").Value = 6
$summary.Range("AG252").Value = 1

